$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new timesheet entry (2014-02-25, 09:45 -> 12:00) is being recorded.
# It belongs right after the last data row (24), before the blank
# separator row and the summary rows, so insert a fresh row at 25 and
# push everything else (blank separator + the 3 summary rows) down.
$ws.Rows("25:25").Insert()

# Populate the new data row.
$ws.Cells.Item(25, 1).Value = 2014
$ws.Cells.Item(25, 2).Value = 2
$ws.Cells.Item(25, 3).Value = 25
$ws.Cells.Item(25, 4).Value = 0.40625
$ws.Cells.Item(25, 5).Value = 0.5
$ws.Cells.Item(25, 6).Formula = "=(E25-D25)*24*60"
$ws.Cells.Item(25, 7).Formula = "=F25/60"

# Reflect the new active selection on the sheet.
$ws.Range("F25").Select()
